$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Host "NOT FOUND: $old"
    }
}

# Rubric total: 60 -> 100 pts
Replace-Text "Rubric (60 pts total)" "Rubric (100 pts total)"

# Title and author: 5 -> 10 pts
Replace-Text "Title and author (5 pts)" "Title and author (10 pts)"

# Introduction: 15 -> 30 pts
Replace-Text "Introduction (15 pts)" "Introduction (30 pts)"

# Description of problem: 5 -> 10 pts
Replace-Text "Description of problem $([char]0x2013) 5 pts" "Description of problem $([char]0x2013) 10 pts"

# Scientific question: 5 -> 10 pts
Replace-Text "Scientific question $([char]0x2013) 5 pts" "Scientific question $([char]0x2013) 10 pts"

# Background info with references: 5 -> 10 pts
Replace-Text "Background information with references $([char]0x2013) 5 pts" "Background information with references $([char]0x2013) 10 pts"

# Methods: 35 -> 50 pts
Replace-Text "Methods (35 pts)" "Methods (50 pts)"

# Background information on the field site: 5 -> 10 pts
Replace-Text "Background information on the field site (5 pts)" "Background information on the field site (10 pts)"

# Map of sampling locations: 5 -> 10 pts
Replace-Text "Map of sampling locations (5 pts)" "Map of sampling locations (10 pts)"

# Field collection: 10 -> 20 pts
Replace-Text "Field collection (10 pts)" "Field collection (20 pts)"

# Description of sample locations: 2.5 -> 5 pts
Replace-Text "Description of sample locations (2.5 pts)" "Description of sample locations (5 pts)"

# Sample collection, including equipment: 2.5 -> 5 pts
Replace-Text "Sample collection, including equipment (2.5 pts)" "Sample collection, including equipment (5 pts)"

# Description of subsampling for different analyses: 2.5 -> 5 pts
Replace-Text "Description of subsampling for different analyses (2.5 pts)" "Description of subsampling for different analyses (5 pts)"

# Field measurements, including equipment: 2.5 -> 5 pts
Replace-Text "Field measurements, including equipment (2.5 pts)" "Field measurements, including equipment (5 pts)"

# Laboratory analysis: 10 -> 5 pts
Replace-Text "Laboratory analysis (10 pts)" "Laboratory analysis (5 pts)"

# Table of samples with field measurements: 5 -> 10 pts
Replace-Text "Table of samples with field measurements (5 pts)" "Table of samples with field measurements (10 pts)"
